$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.609.98"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.75"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.029"
$ws.Range("E4").Value = "  +2.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.12"
$ws.Range("E5").Value = "  +4.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.026"
$ws.Range("E6").Value = "  +2.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4369"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3751"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07401"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8758"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("E11").Value = "  +2.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.860.96"
$ws.Range("E12").Value = "  -3.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.516"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.692"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07178"
$ws.Range("E15").Value = "  +3.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.64"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.031"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.026"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.44"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.592.83"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.254"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.21"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.075.72"
$ws.Range("E24").Value = "  -4.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.34"
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.941"
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.73"
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.308"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.936"
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.17"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09068"
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7695"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.523"
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.875"
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.028"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.155"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01977"
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05287"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.817"
$ws.Range("E40").Value = "  +5.05%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5172"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.731"
$ws.Range("E43").Value = "  +2.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.612"
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "109.00"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.60"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.715"
$ws.Range("E47").Value = "  +3.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4662"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06390"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.883"
$ws.Range("E50").Value = "  +3.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.54"
$ws.Range("E51").Value = "  +5.24%  "
